$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$oldLine1 = "✅ 1000 Bs = 1.73 = 6418.22 pesos"
$newLine1 = "✅ 1000 Bs = 1.76 = 6472.92 pesos"
$oldLine2 = "✅ 6418.22 pesos = 1.73 = 938.6 Bs"
$newLine2 = "✅ 6472.92 pesos = 1.74 = 904.05 Bs"

$text = $wsHoja1.Range("A1").Value()
$text = $text.Replace($oldLine1, $newLine1)
$text = $text.Replace($oldLine2, $newLine2)
$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 568.399
$wsTasas.Range("O10").Value = 3679.2

$wsTasas.Range("N12").Value = 3716.99
$wsTasas.Range("O12").Value = 519.139
